$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value (values that look numeric/percent get a leading
# apostrophe so Excel keeps them as text, matching the source sheet's
# inline-string cells instead of converting them to numbers.
$updates = @(
    @{ Cell = "D2"; Value = '''329.79' }
    @{ Cell = "E2"; Value = '''0.73%' }
    @{ Cell = "G2"; Value = '''12' }
    @{ Cell = "D3"; Value = '''44.30' }
    @{ Cell = "E3"; Value = '''-0.22%' }
    @{ Cell = "G3"; Value = '''12' }
    @{ Cell = "D4"; Value = '''5.500' }
    @{ Cell = "E4"; Value = '''-1.27%' }
    @{ Cell = "G4"; Value = '''12' }
    @{ Cell = "D5"; Value = '''0.08024' }
    @{ Cell = "E5"; Value = '''-0.33%' }
    @{ Cell = "G5"; Value = '''12' }
    @{ Cell = "D6"; Value = '''2.074' }
    @{ Cell = "E6"; Value = '''9.43%' }
    @{ Cell = "G6"; Value = '''12' }
    @{ Cell = "B7"; Value = 'BTSEToken' }
    @{ Cell = "C7"; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = "D7"; Value = '''2.628' }
    @{ Cell = "E7"; Value = '''0.45%' }
    @{ Cell = "G7"; Value = '''12' }
    @{ Cell = "B8"; Value = 'MXToken' }
    @{ Cell = "C8"; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = "D8"; Value = '''0.9539' }
    @{ Cell = "E8"; Value = '''0.85%' }
    @{ Cell = "G8"; Value = '''12' }
    @{ Cell = "D9"; Value = '''0.1144' }
    @{ Cell = "E9"; Value = '''0.28%' }
    @{ Cell = "G9"; Value = '''12' }
    @{ Cell = "D10"; Value = '''0.1894' }
    @{ Cell = "E10"; Value = '''2.74%' }
    @{ Cell = "G10"; Value = '''12' }
    @{ Cell = "D11"; Value = '''10.22' }
    @{ Cell = "E11"; Value = '''6.66%' }
    @{ Cell = "G11"; Value = '''12' }
    @{ Cell = "D12"; Value = '''0.09835' }
    @{ Cell = "E12"; Value = '''1.10%' }
    @{ Cell = "G12"; Value = '''12' }
    @{ Cell = "D13"; Value = '''0.04889' }
    @{ Cell = "E13"; Value = '''12.37%' }
    @{ Cell = "G13"; Value = '''12' }
    @{ Cell = "E14"; Value = '''-0.15%' }
    @{ Cell = "G14"; Value = '''12' }
    @{ Cell = "D15"; Value = '''0.001258' }
    @{ Cell = "E15"; Value = '''-1.10%' }
    @{ Cell = "G15"; Value = '''12' }
    @{ Cell = "D16"; Value = '''0.04093' }
    @{ Cell = "E16"; Value = '''-3.02%' }
    @{ Cell = "G16"; Value = '''12' }
    @{ Cell = "D17"; Value = '''0.006144' }
    @{ Cell = "E17"; Value = '''2.56%' }
    @{ Cell = "G17"; Value = '''12' }
    @{ Cell = "D18"; Value = '''3.389' }
    @{ Cell = "E18"; Value = '''-6.48%' }
    @{ Cell = "G18"; Value = '''12' }
    @{ Cell = "D19"; Value = '''4.404' }
    @{ Cell = "E19"; Value = '''2.15%' }
    @{ Cell = "G19"; Value = '''12' }
    @{ Cell = "D20"; Value = '''0.3395' }
    @{ Cell = "E20"; Value = '''-2.90%' }
    @{ Cell = "G20"; Value = '''12' }
    @{ Cell = "D21"; Value = '''0.1383' }
    @{ Cell = "E21"; Value = '''0.22%' }
    @{ Cell = "G21"; Value = '''12' }
    @{ Cell = "D22"; Value = '''0.2584' }
    @{ Cell = "E22"; Value = '''-2.58%' }
    @{ Cell = "G22"; Value = '''12' }
    @{ Cell = "D23"; Value = '''0.001301' }
    @{ Cell = "E23"; Value = '''4.16%' }
    @{ Cell = "G23"; Value = '''12' }
    @{ Cell = "E24"; Value = '''-2.49%' }
    @{ Cell = "G24"; Value = '''12' }
    @{ Cell = "E25"; Value = '''-4.83%' }
    @{ Cell = "G25"; Value = '''12' }
    @{ Cell = "D26"; Value = '''0.0003747' }
    @{ Cell = "E26"; Value = '''-6.16%' }
    @{ Cell = "G26"; Value = '''12' }
    @{ Cell = "G27"; Value = '''12' }
    @{ Cell = "G28"; Value = '''12' }
    @{ Cell = "G29"; Value = '''12' }
    @{ Cell = "G30"; Value = '''12' }
    @{ Cell = "G31"; Value = '''12' }
    @{ Cell = "G32"; Value = '''12' }
    @{ Cell = "G33"; Value = '''12' }
    @{ Cell = "G34"; Value = '''12' }
    @{ Cell = "G35"; Value = '''12' }
    @{ Cell = "G36"; Value = '''12' }
    @{ Cell = "G37"; Value = '''12' }
    @{ Cell = "D38"; Value = '''0.02589' }
    @{ Cell = "E38"; Value = '''-1.43%' }
    @{ Cell = "G38"; Value = '''12' }
    @{ Cell = "D39"; Value = '''0.05803' }
    @{ Cell = "E39"; Value = '''5.90%' }
    @{ Cell = "G39"; Value = '''12' }
    @{ Cell = "D40"; Value = '''0.007565' }
    @{ Cell = "E40"; Value = '''0.15%' }
    @{ Cell = "G40"; Value = '''12' }
    @{ Cell = "D41"; Value = '''0.1403' }
    @{ Cell = "E41"; Value = '''0.61%' }
    @{ Cell = "G41"; Value = '''12' }
    @{ Cell = "D42"; Value = '''0.007337' }
    @{ Cell = "E42"; Value = '''-0.19%' }
    @{ Cell = "G42"; Value = '''12' }
    @{ Cell = "E43"; Value = '''-0.11%' }
    @{ Cell = "G43"; Value = '''12' }
    @{ Cell = "D44"; Value = '''0.009070' }
    @{ Cell = "E44"; Value = '''2.62%' }
    @{ Cell = "G44"; Value = '''12' }
    @{ Cell = "D45"; Value = '''0.00007037' }
    @{ Cell = "E45"; Value = '''1.46%' }
    @{ Cell = "G45"; Value = '''12' }
    @{ Cell = "E46"; Value = '''-0.07%' }
    @{ Cell = "G46"; Value = '''12' }
    @{ Cell = "D47"; Value = '''0.0005803' }
    @{ Cell = "E47"; Value = '''-0.14%' }
    @{ Cell = "G47"; Value = '''12' }
    @{ Cell = "D48"; Value = '''0.003532' }
    @{ Cell = "E48"; Value = '''55.47%' }
    @{ Cell = "G48"; Value = '''12' }
    @{ Cell = "D49"; Value = '''0.003503' }
    @{ Cell = "E49"; Value = '''-8.61%' }
    @{ Cell = "G49"; Value = '''12' }
    @{ Cell = "E50"; Value = '''-0.07%' }
    @{ Cell = "G50"; Value = '''12' }
    @{ Cell = "E51"; Value = '''-0.07%' }
    @{ Cell = "G51"; Value = '''12' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
